$d = $word.ActiveDocument
$d.Content.Find.Execute("Project :", $true, $false, $false, $false, $false,
                         $true, 1, $false, "Project : ", 2)
